$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: some Price (column D) values look like plain numbers (e.g. "302.99") but
# must stay stored as TEXT, exactly like in the original workbook (they already are
# cells of type string there). Assigning a leading apostrophe (as typed in the Excel
# UI) forces Excel to keep such values as literal text instead of auto-converting
# them to numbers (which would also risk dropping significant trailing zeros, e.g.
# "0.0710" -> 0.071). In a single-quoted PowerShell string, '' represents one
# literal embedded apostrophe, so '''302.99' yields the text  '302.99  i.e. a
# leading apostrophe followed by 302.99.

$ws.Range('D2').Value = '42.605.20'
$ws.Range('E2').Value = '  -2.42%  '

$ws.Range('D3').Value = '2.301.42'
$ws.Range('E3').Value = '  -0.78%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').Value = '''302.99'
$ws.Range('E5').Value = '  -2.62%  '

$ws.Range('D6').Value = '''99.11'
$ws.Range('E6').Value = '  -6.39%  '

$ws.Range('D7').Value = '''0.506'
$ws.Range('E7').Value = '  -5.30%  '

$ws.Range('E8').Value = '  +0.09%  '

$ws.Range('D9').Value = '''0.502'
$ws.Range('E9').Value = '  -5.27%  '

$ws.Range('D10').Value = '''34.67'
$ws.Range('E10').Value = '  -5.94%  '

$ws.Range('E11').Value = '  -3.56%  '

$ws.Range('E12').Value = '  +0.15%  '

$ws.Range('D13').Value = '''6.75'
$ws.Range('E13').Value = '  -4.15%  '

$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.657.92'
$ws.Range('E14').Value = '  -0.70%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '''15.69'
$ws.Range('E15').Value = '  +3.17%  '

$ws.Range('D16').Value = '2.283.34'
$ws.Range('E16').Value = '  -1.71%  '

$ws.Range('E17').Value = '  -1.82%  '

$ws.Range('D18').Value = '42.552.06'
$ws.Range('E18').Value = '  -2.32%  '

$ws.Range('E19').Value = '  -3.46%  '

$ws.Range('D20').Value = '''11.47'
$ws.Range('E20').Value = '  -6.14%  '

$ws.Range('D21').Value = '''6.03'
$ws.Range('E21').Value = '  -2.82%  '

$ws.Range('D22').Value = '''68.06'
$ws.Range('E22').Value = '  -0.61%  '

$ws.Range('D23').Value = '''235.57'
$ws.Range('E23').Value = '  -3.31%  '

$ws.Range('E24').Value = '  -4.03%  '

$ws.Range('D25').Value = '''2.51'
$ws.Range('E25').Value = '  -4.55%  '

$ws.Range('E26').Value = '  +0.18%  '

$ws.Range('D27').Value = '''25.15'
$ws.Range('E27').Value = '  -0.19%  '

$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').Value = '''34.73'
$ws.Range('E28').Value = '  -6.77%  '

$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '''2.07'
$ws.Range('E29').Value = '  -11.99%  '

$ws.Range('E30').Value = '  -5.22%  '

$ws.Range('D31').Value = '''162.39'
$ws.Range('E31').Value = '  -1.92%  '

$ws.Range('D32').Value = '''1.00'
$ws.Range('E32').Value = '  -0.03%  '

$ws.Range('D33').Value = '''5.00'
$ws.Range('E33').Value = '  -5.99%  '

$ws.Range('D34').Value = '''4.61'
$ws.Range('E34').Value = '  +1.28%  '

$ws.Range('E35').Value = '  -5.23%  '

$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D36').Value = '''16.96'
$ws.Range('E36').Value = '  -8.17%  '

$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '''0.0710'
$ws.Range('E37').Value = '  -4.84%  '

$ws.Range('E38').Value = '  -6.18%  '

$ws.Range('E39').Value = '  -5.10%  '

$ws.Range('E40').Value = '  -4.29%  '

$ws.Range('E41').Value = '  -6.88%  '

$ws.Range('D42').Value = '''2.41'
$ws.Range('E42').Value = '  -11.88%  '

$ws.Range('D43').Value = '1.967.72'
$ws.Range('E43').Value = '  -1.90%  '

$ws.Range('E44').Value = '  -4.82%  '

$ws.Range('D45').Value = '''18.52'
$ws.Range('E45').Value = '  -4.20%  '

$ws.Range('D46').Value = '''10.22'
$ws.Range('E46').Value = '  +1.82%  '

$ws.Range('E47').Value = '  -7.58%  '

$ws.Range('D48').Value = '''55.12'
$ws.Range('E48').Value = '  -4.99%  '

$ws.Range('D49').Value = '''2.84'
$ws.Range('E49').Value = '  -2.23%  '

$ws.Range('D50').Value = '2.525.77'
$ws.Range('E50').Value = '  -0.66%  '

$ws.Range('D51').Value = '''4.71'
$ws.Range('E51').Value = '  -0.61%  '
